$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 20 (2024Q4) metrics: total_customers, returning_customers, new_customers
$ws.Range("C20").Value = 270
$ws.Range("D20").Value = 226
$ws.Range("E20").Value = 44

# recurrence_rate = returning_customers / previous_quarter total_customers * 100
$ws.Range("F20").Value = 74.3421052631579
